$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Re-format A13 and A14 so they use the same (duplicate) date style
#    index that the rest of the date column (A2:A12) already uses.
#    We do this by copying the number-format from an already-correct
#    date cell (A12) and pasting only the formatting (not the value)
#    onto A13 / A14.
# ------------------------------------------------------------------
$ws.Cells.Item(12, 1).Copy()
$ws.Cells.Item(13, 1).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(14, 1).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Add the new row (15) for the "Online Stock Span" question that
#    was solved using a monotonic stack.
# ------------------------------------------------------------------

# -- A15: date, formatted like the other date cells in the column.
$ws.Cells.Item(14, 1).Copy()
$ws.Cells.Item(15, 1).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Cells.Item(15, 1).Value = 46079

# -- B15: question title.
$ws.Cells.Item(15, 2).Value = "Online Stock Span"

# -- C15: URL text + hyperlink, styled like the other URL cells.
$ws.Cells.Item(15, 3).Value = "https://leetcode.com/problems/online-stock-span/"
$ws.Hyperlinks.Add($ws.Cells.Item(15, 3), "https://leetcode.com/problems/online-stock-span/")

# Hyperlinks.Add() forces Excel's built-in "Hyperlink" cell style onto
# the target cell/style table. Re-apply the plain formatting used by
# the rest of column C (same underline/theme font, no named style) and
# drop the now-unused "Hyperlink" named style it introduced so the
# cell matches its neighbours (C2:C14).
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$wb.Styles.Item("Hyperlink").Delete()

Write-Output "Applied Online Stock Span row"
